# Re-pull / push updated dSF (column F) values for a set of rows.
# Mirrors: "repull data, push all data, mean calculation"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of worksheet row number -> new value for column F (dSF)
$updates = @{
    5  = 1
    7  = 3
    13 = 2
    18 = 3
    20 = 3
    23 = 3
    31 = -1
    33 = 1
    35 = -2
    38 = 1
    40 = 2
    41 = 3
    52 = -5
    56 = 0
    59 = -2
    60 = -1
    67 = -8
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
